$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 31   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/11/2024  Through  11/17/2024"

# --- Donor cells used to copy styles without creating new style entries ---
# s=14 (#,##0 integer) donor: I14   |  s=15 (#,##0.0 decimal) donor: K14
# s=13 (General/text) + shared-string "0" donor: C14   |  "***.* " donor: E14
$fmt14 = $ws.Range("I14").NumberFormat
$fmt15 = $ws.Range("K14").NumberFormat

# --- Cells that flip from a numeric style to the text "N/A" placeholders ---
# (style must become s=13 + shared string, so we copy a cell that already has that exact combo)
$ws.Range("C14").Copy($ws.Range("C23"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))

# --- Cells that flip from a text "N/A" placeholder to an actual number ---
# (set the destination NumberFormat to the donors format string first so the engine reuses
#  the existing numeric style instead of fabricating a new one, then assign the number)
$ws.Range("D18").NumberFormat = $fmt14
$ws.Range("D18").Value = 1
$ws.Range("E18").NumberFormat = $fmt15
$ws.Range("E18").Value = 0
$ws.Range("D23").NumberFormat = $fmt14
$ws.Range("D23").Value = 1
$ws.Range("E23").NumberFormat = $fmt15
$ws.Range("E23").Value = -100
$ws.Range("C28").NumberFormat = $fmt14
$ws.Range("C28").Value = 1
$ws.Range("C33").NumberFormat = $fmt14
$ws.Range("C33").Value = 1
$ws.Range("F33").NumberFormat = $fmt14
$ws.Range("F33").Value = 1
$ws.Range("I33").NumberFormat = $fmt14
$ws.Range("I33").Value = 1

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("N15").Value = -65.217391304347
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 66.666666666666
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 200
$ws.Range("I16").Value = 112
$ws.Range("J16").Value = 96
$ws.Range("K16").Value = 16.666666666666
$ws.Range("L16").Value = 10.891089108910
$ws.Range("M16").Value = -42.857142857142
$ws.Range("N16").Value = -81.788617886178
$ws.Range("D17").Value = 2
$ws.Range("F17").Value = 8
$ws.Range("H17").Value = 14.285714285714
$ws.Range("J17").Value = 142
$ws.Range("K17").Value = 8.450704225352
$ws.Range("L17").Value = 2.666666666666
$ws.Range("M17").Value = 35.087719298245
$ws.Range("N17").Value = -45.964912280701
$ws.Range("C18").Value = 1
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 300
$ws.Range("I18").Value = 66
$ws.Range("J18").Value = 80
$ws.Range("K18").Value = -17.5
$ws.Range("L18").Value = -37.142857142857
$ws.Range("M18").Value = -73.6
$ws.Range("N18").Value = -93.542074363992
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = -23.728813559322
$ws.Range("I19").Value = 514
$ws.Range("J19").Value = 587
$ws.Range("K19").Value = -12.436115843270
$ws.Range("L19").Value = 11.255411255411
$ws.Range("M19").Value = 15.246636771300
$ws.Range("N19").Value = -9.347442680776
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -21.428571428571
$ws.Range("I20").Value = 157
$ws.Range("J20").Value = 118
$ws.Range("K20").Value = 33.050847457627
$ws.Range("L20").Value = 36.521739130434
$ws.Range("M20").Value = 9.790209790209
$ws.Range("N20").Value = -93.533772652388
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -13.636363636363
$ws.Range("F21").Value = 87
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = -2.247191011235
$ws.Range("I21").Value = 1013
$ws.Range("J21").Value = 1039
$ws.Range("K21").Value = -2.502406159769
$ws.Range("L21").Value = 7.537154989384
$ws.Range("M21").Value = -13.196229648671
$ws.Range("N21").Value = -79.551877270892
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -50
$ws.Range("J23").Value = 27
$ws.Range("K23").Value = 0
$ws.Range("C24").Value = 38
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 26.666666666666
$ws.Range("F24").Value = 117
$ws.Range("G24").Value = 116
$ws.Range("H24").Value = 0.862068965517
$ws.Range("I24").Value = 1449
$ws.Range("J24").Value = 1160
$ws.Range("K24").Value = 24.913793103448
$ws.Range("L24").Value = 43.892750744786
$ws.Range("M24").Value = 61.71875
$ws.Range("C25").Value = 31
$ws.Range("D25").Value = 27
$ws.Range("E25").Value = 14.814814814814
$ws.Range("F25").Value = 93
$ws.Range("G25").Value = 106
$ws.Range("H25").Value = -12.264150943396
$ws.Range("I25").Value = 1200
$ws.Range("J25").Value = 858
$ws.Range("K25").Value = 39.860139860139
$ws.Range("L25").Value = 87.5
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 25
$ws.Range("F26").Value = 22
$ws.Range("G26").Value = 27
$ws.Range("H26").Value = -18.518518518518
$ws.Range("I26").Value = 275
$ws.Range("J26").Value = 241
$ws.Range("K26").Value = 14.107883817427
$ws.Range("L26").Value = 30.952380952381
$ws.Range("M26").Value = -13.249211356466
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 33.333333333333
$ws.Range("I28").Value = 37
$ws.Range("K28").Value = 37.037037037037
$ws.Range("L28").Value = 27.586206896551
$ws.Range("I31").Value = 9
$ws.Range("K31").Value = 350
$ws.Range("L31").Value = 50
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = -66.666666666666

Write-Host "edit.ps1 completed"
